# Apply invoice edits: add a client discount amount and update the
# footer contact placeholder with the RPA developer's name/email.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice")

# Set the "Client discount" unit price (row 18, column E) so the
# discount line contributes to the subtotal/tax/total calculations.
$ws.Range("E18").Value = 100

# Update the footer placeholder text with the RPA developer contact info.
$ws.Range("A31").Value = "RPA Developer - RPADeveloper@Uipath.com"
